# Generate Report for Handoff
#
# Inserts a new tracked file, "86297059-e3e2-48e1-a692-ce334266629a.md",
# into the localization-status report. The new entry is placed immediately
# before the existing "bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.md" row (and
# before the trailing ".localization-config" row) on all three sheets:
# Overview, zh-cn and de-de.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview  (File Name | zh-cn | de-de)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Rows.Item(7).Insert()

$ws1.Range("A7").Value2 = "86297059-e3e2-48e1-a692-ce334266629a.md"
$ws1.Range("B7").Value2 = "Ready for handoff"
$ws1.Range("C7").Value2 = "Ready for handoff"

# Rebuild the hyperlinks for this sheet in the correct final order (the
# simulated Hyperlinks collection only supports whole-sheet clear + append).
$ws1.Range("A1:C20").Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/979c623da3952ace064a5e4caf6972e7d2549fa3/e2e/27f7f665-9af6-4ffe-a2bb-707830a84b7a.md", "", "", "27f7f665-9af6-4ffe-a2bb-707830a84b7a.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0f1961626b414a1cddc99ed1e65daaf28ba3bf81/e2e/31935e31-00e9-473d-9e54-6e79352372b7.md", "", "", "31935e31-00e9-473d-9e54-6e79352372b7.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0f1961626b414a1cddc99ed1e65daaf28ba3bf81/e2e/adb73576-fdd6-49ca-96f2-83f54a1e4446.md", "", "", "adb73576-fdd6-49ca-96f2-83f54a1e4446.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7ab5f779eab5a2c40381f5f595326f0a0a2458d4/e2e/edb7b64d-b387-4334-ae49-3b80715bbcc2.md", "", "", "edb7b64d-b387-4334-ae49-3b80715bbcc2.md")
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/84435e29e6a354b3f815aee67b2e6ea41b6596b8/e2e/4cbda6e1-9396-404e-ae9a-df7f4d1ca222.md", "", "", "4cbda6e1-9396-404e-ae9a-df7f4d1ca222.md")
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/383297a008d98f3970af6ee606767008952a8124/e2e/86297059-e3e2-48e1-a692-ce334266629a.md", "", "", "86297059-e3e2-48e1-a692-ce334266629a.md")
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/8267d1c098713c723f9cc24cd0aaa33c5fac1b4d/e2e/bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.md", "", "", "bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.md")
$ws1.Hyperlinks.Add($ws1.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/84435e29e6a354b3f815aee67b2e6ea41b6596b8/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(7).Insert()

$ws2.Range("A7").Value2 = "86297059-e3e2-48e1-a692-ce334266629a.md"
$ws2.Range("B7").Value2 = "Ready for handoff"
$ws2.Range("C7").Value2 = "86297059-e3e2-48e1-a692-ce334266629a.02e5e4b54511cb3b3766efa871fce9d1dc61990a.zh-cn.xlf"
$ws2.Range("D7").Value2 = "2016-03-09 15:15:29"
$ws2.Range("G7").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H7").Value2 = "Include"

$ws2.Range("A1:I20").Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/979c623da3952ace064a5e4caf6972e7d2549fa3/e2e/27f7f665-9af6-4ffe-a2bb-707830a84b7a.md", "", "", "27f7f665-9af6-4ffe-a2bb-707830a84b7a.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c2ec876f824e813ba4283aeebea261b84fc8b0b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/27f7f665-9af6-4ffe-a2bb-707830a84b7a.764d3d58d9e37c89a16951bf3ce81701003ec11f.zh-cn.xlf", "", "", "27f7f665-9af6-4ffe-a2bb-707830a84b7a.764d3d58d9e37c89a16951bf3ce81701003ec11f.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/40eefd2062ad86bd3781fe442b2ef7badaf264d3/e2e/27f7f665-9af6-4ffe-a2bb-707830a84b7a.md", "", "", "27f7f665-9af6-4ffe-a2bb-707830a84b7a.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f4a1c21066bc9e9548c5879d5ae45111c52e027c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/27f7f665-9af6-4ffe-a2bb-707830a84b7a.764d3d58d9e37c89a16951bf3ce81701003ec11f.zh-cn.xlf", "", "", "27f7f665-9af6-4ffe-a2bb-707830a84b7a.764d3d58d9e37c89a16951bf3ce81701003ec11f.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0f1961626b414a1cddc99ed1e65daaf28ba3bf81/e2e/31935e31-00e9-473d-9e54-6e79352372b7.md", "", "", "31935e31-00e9-473d-9e54-6e79352372b7.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e502f16a2f22ed8587fd4cecdfc1a3dbe80889a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/31935e31-00e9-473d-9e54-6e79352372b7.64e33c10c33d0819d892a85390aee23da1040aa3.zh-cn.xlf", "", "", "31935e31-00e9-473d-9e54-6e79352372b7.64e33c10c33d0819d892a85390aee23da1040aa3.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0f1961626b414a1cddc99ed1e65daaf28ba3bf81/e2e/adb73576-fdd6-49ca-96f2-83f54a1e4446.md", "", "", "adb73576-fdd6-49ca-96f2-83f54a1e4446.md")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e502f16a2f22ed8587fd4cecdfc1a3dbe80889a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/adb73576-fdd6-49ca-96f2-83f54a1e4446.124d8101d6ff29239d57603d2b23717faab74709.zh-cn.xlf", "", "", "adb73576-fdd6-49ca-96f2-83f54a1e4446.124d8101d6ff29239d57603d2b23717faab74709.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7ab5f779eab5a2c40381f5f595326f0a0a2458d4/e2e/edb7b64d-b387-4334-ae49-3b80715bbcc2.md", "", "", "edb7b64d-b387-4334-ae49-3b80715bbcc2.md")
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/631a28965a08bbb19ea65d0f661eb686792b41e5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/edb7b64d-b387-4334-ae49-3b80715bbcc2.feee2006d607e5d0c4a04728cfa67b44ad4c2842.zh-cn.xlf", "", "", "edb7b64d-b387-4334-ae49-3b80715bbcc2.feee2006d607e5d0c4a04728cfa67b44ad4c2842.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/84435e29e6a354b3f815aee67b2e6ea41b6596b8/e2e/4cbda6e1-9396-404e-ae9a-df7f4d1ca222.md", "", "", "4cbda6e1-9396-404e-ae9a-df7f4d1ca222.md")
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9fc5fadee8e2f2fff7ad903bdcaf9774530318c9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4cbda6e1-9396-404e-ae9a-df7f4d1ca222.fdae4a8ba869d2b46b3d3714d4cb9d8a2a763620.zh-cn.xlf", "", "", "4cbda6e1-9396-404e-ae9a-df7f4d1ca222.fdae4a8ba869d2b46b3d3714d4cb9d8a2a763620.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/383297a008d98f3970af6ee606767008952a8124/e2e/86297059-e3e2-48e1-a692-ce334266629a.md", "", "", "86297059-e3e2-48e1-a692-ce334266629a.md")
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4d410487f0df48a27d9c18e13ee391b6a854f10e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/86297059-e3e2-48e1-a692-ce334266629a.02e5e4b54511cb3b3766efa871fce9d1dc61990a.zh-cn.xlf", "", "", "86297059-e3e2-48e1-a692-ce334266629a.02e5e4b54511cb3b3766efa871fce9d1dc61990a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/8267d1c098713c723f9cc24cd0aaa33c5fac1b4d/e2e/bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.md", "", "", "bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.md")
$ws2.Hyperlinks.Add($ws2.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/93f3aa66410f2ae2936759f41c3b95689a3869b4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.38cde8c684e5a7400fe076009a976c113630a202.zh-cn.xlf", "", "", "bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.38cde8c684e5a7400fe076009a976c113630a202.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/84435e29e6a354b3f815aee67b2e6ea41b6596b8/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(7).Insert()

$ws3.Range("A7").Value2 = "86297059-e3e2-48e1-a692-ce334266629a.md"
$ws3.Range("B7").Value2 = "Ready for handoff"
$ws3.Range("C7").Value2 = "86297059-e3e2-48e1-a692-ce334266629a.02e5e4b54511cb3b3766efa871fce9d1dc61990a.de-de.xlf"
$ws3.Range("D7").Value2 = "2016-03-09 15:15:40"
$ws3.Range("G7").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H7").Value2 = "Include"

$ws3.Range("A1:I20").Hyperlinks.Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/979c623da3952ace064a5e4caf6972e7d2549fa3/e2e/27f7f665-9af6-4ffe-a2bb-707830a84b7a.md", "", "", "27f7f665-9af6-4ffe-a2bb-707830a84b7a.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed37e8ea585d292a20e47faa3c3f7d46e9505114/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/27f7f665-9af6-4ffe-a2bb-707830a84b7a.764d3d58d9e37c89a16951bf3ce81701003ec11f.de-de.xlf", "", "", "27f7f665-9af6-4ffe-a2bb-707830a84b7a.764d3d58d9e37c89a16951bf3ce81701003ec11f.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/96ef2bb2fabf4c56fe3fe80a1e207e4326033d30/e2e/27f7f665-9af6-4ffe-a2bb-707830a84b7a.md", "", "", "27f7f665-9af6-4ffe-a2bb-707830a84b7a.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/be4ba82fe78edec995c5ab97bfc636bb224af6ab/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/27f7f665-9af6-4ffe-a2bb-707830a84b7a.764d3d58d9e37c89a16951bf3ce81701003ec11f.de-de.xlf", "", "", "27f7f665-9af6-4ffe-a2bb-707830a84b7a.764d3d58d9e37c89a16951bf3ce81701003ec11f.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0f1961626b414a1cddc99ed1e65daaf28ba3bf81/e2e/31935e31-00e9-473d-9e54-6e79352372b7.md", "", "", "31935e31-00e9-473d-9e54-6e79352372b7.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09f98c4762c3e21c2f15a8235d48a95572f705d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/31935e31-00e9-473d-9e54-6e79352372b7.64e33c10c33d0819d892a85390aee23da1040aa3.de-de.xlf", "", "", "31935e31-00e9-473d-9e54-6e79352372b7.64e33c10c33d0819d892a85390aee23da1040aa3.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0f1961626b414a1cddc99ed1e65daaf28ba3bf81/e2e/adb73576-fdd6-49ca-96f2-83f54a1e4446.md", "", "", "adb73576-fdd6-49ca-96f2-83f54a1e4446.md")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/09f98c4762c3e21c2f15a8235d48a95572f705d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/adb73576-fdd6-49ca-96f2-83f54a1e4446.124d8101d6ff29239d57603d2b23717faab74709.de-de.xlf", "", "", "adb73576-fdd6-49ca-96f2-83f54a1e4446.124d8101d6ff29239d57603d2b23717faab74709.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7ab5f779eab5a2c40381f5f595326f0a0a2458d4/e2e/edb7b64d-b387-4334-ae49-3b80715bbcc2.md", "", "", "edb7b64d-b387-4334-ae49-3b80715bbcc2.md")
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/88fa522ae3dff4aa1ca1bf876580a09ba0c75f00/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/edb7b64d-b387-4334-ae49-3b80715bbcc2.feee2006d607e5d0c4a04728cfa67b44ad4c2842.de-de.xlf", "", "", "edb7b64d-b387-4334-ae49-3b80715bbcc2.feee2006d607e5d0c4a04728cfa67b44ad4c2842.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/84435e29e6a354b3f815aee67b2e6ea41b6596b8/e2e/4cbda6e1-9396-404e-ae9a-df7f4d1ca222.md", "", "", "4cbda6e1-9396-404e-ae9a-df7f4d1ca222.md")
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/032a1bdcaf6e877dc46bd3b438f51daf08d7bf16/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4cbda6e1-9396-404e-ae9a-df7f4d1ca222.fdae4a8ba869d2b46b3d3714d4cb9d8a2a763620.de-de.xlf", "", "", "4cbda6e1-9396-404e-ae9a-df7f4d1ca222.fdae4a8ba869d2b46b3d3714d4cb9d8a2a763620.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/383297a008d98f3970af6ee606767008952a8124/e2e/86297059-e3e2-48e1-a692-ce334266629a.md", "", "", "86297059-e3e2-48e1-a692-ce334266629a.md")
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f683870131309f32e499216a3376b60c0bf8021b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/86297059-e3e2-48e1-a692-ce334266629a.02e5e4b54511cb3b3766efa871fce9d1dc61990a.de-de.xlf", "", "", "86297059-e3e2-48e1-a692-ce334266629a.02e5e4b54511cb3b3766efa871fce9d1dc61990a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/8267d1c098713c723f9cc24cd0aaa33c5fac1b4d/e2e/bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.md", "", "", "bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.md")
$ws3.Hyperlinks.Add($ws3.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d2ed0e236dc19e430293fff3af29fb43382515a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.38cde8c684e5a7400fe076009a976c113630a202.de-de.xlf", "", "", "bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b.38cde8c684e5a7400fe076009a976c113630a202.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/84435e29e6a354b3f815aee67b2e6ea41b6596b8/.localization-config", "", "", ".localization-config")

Write-Output "Inserted 86297059-e3e2-48e1-a692-ce334266629a.md handoff rows on Overview, zh-cn and de-de sheets."
